# Repull data: update the dSF (column F) values for the rows whose
# figures changed after re-pulling the source data / recalculating
# the mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -6
    7  = -2
    9  = -5
    10 = -4
    12 = -1
    14 = 0
    24 = -8
    26 = -1
    29 = -7
    30 = 0
    32 = -3
    33 = 3
    35 = 5
    37 = 3
    42 = 0
    45 = 5
    53 = -5
    55 = 3
    56 = -5
    59 = 2
    60 = -2
    61 = 3
    62 = 2
    64 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value2 = $updates[$row]
}
